$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows down
$ws.Rows.Item(21).Insert()

# Fill in the new row 21 with data
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 44811
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100102
$ws.Cells.Item(21, 8).Value = "Cítricos"
$ws.Cells.Item(21, 9).Value = 100102004
$ws.Cells.Item(21, 10).Value = "Mandarina"
$ws.Cells.Item(21, 11).Value = "Murcott"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 250
$ws.Cells.Item(21, 14).Value = 16000
$ws.Cells.Item(21, 15).Value = 17000
$ws.Cells.Item(21, 16).Value = 16500
$ws.Cells.Item(21, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(21, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(21, 19).Value = 825
$ws.Cells.Item(21, 20).Value = 20

Write-Host "done"
